# Add new columns O:R (Matte, Gloss, Soft Touch, Varnish) with data to Sheet1,
# and make Sheet1 the selected/active tab (previously Sheet3 was active).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# Apply the style to O2 first (numFmt 165 "Comma, 0dp"; fontId 0; fill theme5
# tint -0.25; border: medium on the left edge only). Reuse the existing format
# set from J2 (same numFmt/font/fill, different border) via Copy/PasteSpecial
# so the shared numFmt/fill/font table entries are reused rather than
# duplicated, then adjust just the border to a left-only medium edge.
$ws1.Range("J2").Copy()
$ws1.Range("O2").PasteSpecial(-4122)
$ws1.Range("O2").Borders.LineStyle = -4142
$ws1.Range("O2").Borders(7).Weight = -4138

# New header values in row 1, columns O:R
$ws1.Range("O1").Value = "Matte"
$ws1.Range("P1").Value = "Gloss"
$ws1.Range("Q1").Value = "Soft Touch"
$ws1.Range("R1").Value = "Varnish"

# New data values in row 2, columns O:R (set after the paste/style step above
# so these values aren't clobbered by pasting J2's content)
$ws1.Range("O2").Value = 3.5
$ws1.Range("P2").Value = 2.5
$ws1.Range("Q2").Value = 25
$ws1.Range("R2").Value = 6

# Update selection on Sheet1
$ws1.Range("S2").Select()

# Make Sheet1 the active/selected sheet (was Sheet3)
$ws1.Activate()
